$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the date value stored in A101 (existing row) ---
# (previous value 45482.6292013889 -> new value 45482.2916666667)
$ws.Range("A101").Value = 45482.2916666667

# --- Append a new data row (row 102) with the latest OHLCV results ---
$ws.Range("A102").Value = 45483.6241550926
$ws.Range("B102").Value = 6000
$ws.Range("C102").Value = 6.05999994277954
$ws.Range("D102").Value = 6
$ws.Range("E102").Value = 6.03999996185303
$ws.Range("F102").Value = 6
# G (adj_close) and H (ticker) are stored as text in this workbook, so force
# text entry the same way a real user would (leading apostrophe) rather than
# letting Excel auto-detect "6" as a number.
$ws.Range("G102").Value = "'6"
$ws.Range("G102").Style = "Normal"
$ws.Range("H102").Value = "PAL.MI"

# Reuse A101's existing date-time style/format for the new date cell instead
# of creating a brand new style entry.
$ws.Range("A101").Copy() | Out-Null
$ws.Range("A102").PasteSpecial(-4122) | Out-Null

# Clear the clipboard marquee left behind by Copy/PasteSpecial.
$excel.CutCopyMode = 0
